# Adds the "ODI Bowling Extra" worksheet (MATCH_CODE / MAIDEN_OVERS /
# PERCENT_WICKETS_OF_ALL) as the new last sheet in the workbook, mirroring
# the existing "ODI Batting Extra" sheet layout, and populates it with the
# scraped bowling data for each MATCH_CODE already present in the workbook.

$wb = $excel.ActiveWorkbook

# Use "ODI Batting Extra" (the existing "Extra" sheet) both as the anchor to
# insert after, and as the source of header formatting to copy.
$sourceSheet = $wb.Worksheets.Item("ODI Batting Extra")

$newSheet = $wb.Worksheets.Add($null, $sourceSheet)
$newSheet.Name = "ODI Bowling Extra"

# --- Header row -------------------------------------------------------
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "MAIDEN_OVERS"
$newSheet.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Match the bold/bordered header styling already used on the other sheets.
$sourceSheet.Range("A1:C1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
# Keep MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL as plain text
# (matching the rest of the scraped workbook, which stores these scraped
# fields as text rather than numbers/percentages) by forcing the range to
# Text format before any values are written into it.
$dataRange = $newSheet.Range("A2:C21")
$dataRange.NumberFormat = "@"

$rows = @(
    @("4573", "",  ""),
    @("4575", "0", ""),
    @("4576", "0", "20.00%"),
    @("4581", "1", "20.00%"),
    @("4604", "0", "50.00%"),
    @("4610", "0", ""),
    @("4612", "2", ""),
    @("4617", "2", ""),
    @("4625", "0", ""),
    @("4629", "0", "10.00%"),
    @("4631", "", ""),
    @("4632", "1", "10.00%"),
    @("4635", "", ""),
    @("4677", "1", "30.00%"),
    @("4681", "1", "10.00%"),
    @("4680", "0", "10.00%"),
    @("4684", "0", "20.00%"),
    @("4702", "2", "30.00%"),
    @("4705", "0", ""),
    @("4706", "", "")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $newSheet.Cells.Item($r, 1).Value = $rows[$i][0]
    if ($rows[$i][1] -ne "") {
        $newSheet.Cells.Item($r, 2).Value = $rows[$i][1]
    }
    if ($rows[$i][2] -ne "") {
        $newSheet.Cells.Item($r, 3).Value = $rows[$i][2]
    }
}

$newSheet.Range("A1").Select()
